$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style index 1 = "Hyperlink" cell style (already used by B2:B5, B11)
# Style index 2 = new plain style with applyAlignment - will be created by setting
# a horizontal alignment on E7 (new xf entry appended to cellXfs)

# Row 6
$ws.Range("B6").Value = "jisola.tsoft@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("C6").Value = 12061990
$ws.Range("D6").Value = "juan martin isola"

# Row 7
$ws.Range("B7").Value = "jisola.tsoft@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("C7").Value = 12061990
$ws.Range("D7").Value = "juan martin isola"
$ws.Range("E7").Value = "Cancelar solicitud"
$ws.Range("E7").HorizontalAlignment = -4108

# Row 8
$ws.Range("B8").Value = "jisola.tsoft@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("C8").Value = 12061990
$ws.Range("D8").Value = "juan martin isola"
$ws.Range("E8").Value = "Agregar"

# Row 9
$ws.Range("B9").Value = "jisola.tsoft@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B9").Style = "Hyperlink"
$ws.Range("C9").Value = 12061990
$ws.Range("D9").Value = "Tsoft"
$ws.Range("E9").Value = "Te gusta"

# Update selection to F6 (cosmetic, matches diff)
$ws.Range("F6").Select() | Out-Null

# Rename hyperlink cell style from "Hyperlink" to "Hipervínculo"
$wb.Styles.Item("Hyperlink").Name = "Hipervínculo"
